$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.032808780670166
$ws.Range("B1").Value = 2.296585559844971
$ws.Range("C1").Value = 2.312442302703857
$ws.Range("D1").Value = 2.767592191696167
$ws.Range("E1").Value = 1.260187387466431
